# Update dSF (column F) values for the indicated rows on Sheet1.
# This reflects a "repull data, push all data, mean calculation" update
# where only the dSF column changed relative to the previously saved data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7  = -2
    8  = -7
    9  = -5
    12 = -5
    14 = 0
    15 = -3
    17 = -1
    19 = 6
    23 = -4
    24 = -10
    27 = -1
    32 = -14
    34 = 1
    39 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
